# Apply the edit described by the diff: on the "АПРЕЛЬ 2025" sheet, the
# contents of columns E:G (for every data row) are rotated one slot to the
# right -- the value that used to sit in G (the "08.06." date-stamp) moves
# into E, the old E value moves into F, and the old F value moves into G.
# Rows 1 and 8 are merged header/date rows and are not touched; column D
# (birth dates) is untouched as well. The "МАЙ 2025" sheet has no cell
# content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("АПРЕЛЬ 2025")

$dataRows = @(2, 3, 4, 5, 6, 7, 9, 10, 11, 12, 13, 14)

foreach ($r in $dataRows) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $gCell = $ws.Cells.Item($r, 7)   # column G

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2
    $gVal = $gCell.Value2

    # rotate right: E <- G, F <- E, G <- F
    $eCell.Value2 = $gVal
    $fCell.Value2 = $eVal
    $gCell.Value2 = $fVal
}

# Reflect the author's final selection state (whole sheet selected, active
# cell left at J13) on the "АПРЕЛЬ 2025" sheet.
$ws.Range("A1:XFD1048576").Select() | Out-Null
